# Updates the cryptocurrency price/volume table on Sheet1 with freshly
# scraped values (GitHub Actions data refresh), and fixes the ordering of
# the RenderToken / PaxDollar rows (46-47) whose Coin name, Link, Price and
# Volume values had been swapped.
#
# Numeric-looking "Price" values (column D) are explicitly written as text
# (NumberFormat "@") so Excel does not silently reinterpret strings such as
# "1.001" as a floating point number, then the format is reset back to the
# default "Normal" style so no stray cell styling is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.830.14'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').Value = '1.891.51'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7743'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3144'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07479'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.40'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08122'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7671'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.455'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').Value = '1.907.16'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.206'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '29.862.87'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007884'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.102'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('D23').Value = '2.118.55'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1578'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.431'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.81'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.044'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.434'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.550'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.491'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.097'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05516'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.254'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7567'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01922'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').Value = '1.164.24'
$ws.Range('E41').Value = '  +12.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4453'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '73.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.960'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8473'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.901'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.937'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.085'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.551'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.22%  '
